$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''65.233.73'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -2.48%  '
$ws.Range("D3").Value = '''3.373.61'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -2.78%  '
$ws.Range("E4").Value = '  +0.18%  '
$ws.Range("D5").Value = '''592.64'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.84%  '
$ws.Range("D6").Value = '''140.57'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -5.29%  '
$ws.Range("D7").Value = '''0.999'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.06%  '
$ws.Range("D8").Value = '''3.373.39'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.78%  '
$ws.Range("D9").Value = '''0.467'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -3.44%  '
$ws.Range("D10").Value = '''0.133'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -6.72%  '
$ws.Range("D11").Value = '''7.87'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.88%  '
$ws.Range("D12").Value = '''0.404'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -4.89%  '
$ws.Range("D13").Value = '''3.962.46'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.30%  '
$ws.Range("B14").Value = 'Avalanche'
$ws.Range("C14").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D14").Value = '''29.64'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -6.94%  '
$ws.Range("B15").Value = 'ShibaInu'
$ws.Range("C15").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D15").Value = '''0.0000198'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -8.33%  '
$ws.Range("E16").Value = '  -0.48%  '
$ws.Range("B17").Value = 'WrappedBTC'
$ws.Range("C17").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D17").Value = '''65.042.04'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.69%  '
$ws.Range("B18").Value = 'WrappedEther'
$ws.Range("C18").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D18").Value = '''3.383.18'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.37%  '
$ws.Range("D19").Value = '''10.29'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.36%  '
$ws.Range("D20").Value = '''6.07'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -6.09%  '
$ws.Range("D21").Value = '''14.60'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -5.28%  '
$ws.Range("D22").Value = '''414.30'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -5.80%  '
$ws.Range("D23").Value = '''0.577'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -5.90%  '
$ws.Range("D24").Value = '''77.14'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.29%  '
$ws.Range("D26").Value = '''3.511.39'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.62%  '
$ws.Range("D27").Value = '''0.0000108'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -10.89%  '
$ws.Range("D28").Value = '''9.18'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -6.52%  '
$ws.Range("D29").Value = '''7.72'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -8.95%  '
$ws.Range("D30").Value = '''2.41'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.00%  '
$ws.Range("E31").Value = '  -0.10%  '
$ws.Range("E32").Value = '  -4.64%  '
$ws.Range("D33").Value = '''1.45'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -10.32%  '
$ws.Range("B34").Value = 'RenzoRestakedETH'
$ws.Range("C34").Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range("D34").Value = '''3.377.02'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.41%  '
$ws.Range("B35").Value = 'EthereumClassic'
$ws.Range("C35").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D35").Value = '''24.19'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.98%  '
$ws.Range("E36").Value = '  -0.05%  '
$ws.Range("D37").Value = '''1.67'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -7.85%  '
$ws.Range("D38").Value = '''5.47'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -10.24%  '
$ws.Range("D39").Value = '''7.50'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -5.55%  '
$ws.Range("E40").Value = '  +0.30%  '
$ws.Range("D41").Value = '''165.68'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -4.59%  '
$ws.Range("D42").Value = '''0.0851'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -4.77%  '
$ws.Range("B43").Value = 'Mantle'
$ws.Range("C43").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D43").Value = '''0.864'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.64%  '
$ws.Range("B44").Value = 'Filecoin'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D44").Value = '''5.00'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -7.74%  '
$ws.Range("D45").Value = '''1.90'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -12.42%  '
$ws.Range("D46").Value = '''45.27'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.91%  '
$ws.Range("D47").Value = '''26.21'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -11.11%  '
$ws.Range("D48").Value = '''1.16'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -7.17%  '
$ws.Range("D49").Value = '''7.03'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -5.95%  '
$ws.Range("D50").Value = '''2.24'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -9.94%  '
$ws.Range("B51").Value = 'TheGraph'
$ws.Range("C51").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D51").Value = '''0.232'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -6.35%  '
